# Update "想去人数" (interest count) values in column F across sheets
# to match the freshly regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 13503
$ws1.Cells.Item(5, 6).Value = 788
$ws1.Cells.Item(13, 6).Value = 24506
$ws1.Cells.Item(20, 6).Value = 327
$ws1.Cells.Item(24, 6).Value = 246
$ws1.Cells.Item(25, 6).Value = 295
$ws1.Cells.Item(27, 6).Value = 1378
$ws1.Cells.Item(28, 6).Value = 88

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 4487
$ws2.Cells.Item(6, 6).Value = 28

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 4604
$ws3.Cells.Item(4, 6).Value = 125

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 13503
$ws4.Cells.Item(6, 6).Value = 788
$ws4.Cells.Item(7, 6).Value = 4604
$ws4.Cells.Item(14, 6).Value = 125
$ws4.Cells.Item(15, 6).Value = 24506
$ws4.Cells.Item(17, 6).Value = 4487
$ws4.Cells.Item(25, 6).Value = 28
$ws4.Cells.Item(32, 6).Value = 327
$ws4.Cells.Item(37, 6).Value = 246
$ws4.Cells.Item(40, 6).Value = 295
$ws4.Cells.Item(43, 6).Value = 1378
$ws4.Cells.Item(44, 6).Value = 88
